$wb = $excel.ActiveWorkbook

# Grupo A
$ws = $wb.Worksheets.Item("Grupo A")
$ws.Range("B2").Value = "Tatols Beants F.C"
$ws.Range("B3").Value = "JV5 Tricolor Gaúcho"
$ws.Range("B4").Value = "JUV. KP"
$ws.Range("B5").Value = "SERGRILLO"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 64.95999999999999
$ws.Range("H2").Value = 41.6
$ws.Range("I2").Value = 23.35999999999999
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 39.66
$ws.Range("H3").Value = 38.5
$ws.Range("I3").Value = 1.159999999999997
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 41.6
$ws.Range("H4").Value = 64.95999999999999
$ws.Range("I4").Value = -23.35999999999999
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 38.5
$ws.Range("H5").Value = 39.66
$ws.Range("I5").Value = -1.159999999999997

# Grupo B
$ws = $wb.Worksheets.Item("Grupo B")
$ws.Range("B2").Value = "S.E.R. GRILLO"
$ws.Range("B3").Value = "Dom Camillo68"
$ws.Range("B4").Value = "Máquina Laranjja"
$ws.Range("B5").Value = "LISI GREMISTA"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 59.76
$ws.Range("H2").Value = 45.86
$ws.Range("I2").Value = 13.9
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 58.4
$ws.Range("H3").Value = 39.66
$ws.Range("I3").Value = 18.74
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 45.86
$ws.Range("H4").Value = 59.76
$ws.Range("I4").Value = -13.9
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 39.66
$ws.Range("H5").Value = 58.4
$ws.Range("I5").Value = -18.74

# Grupo C
$ws = $wb.Worksheets.Item("Grupo C")
$ws.Range("B2").Value = "dasdoresfc"
$ws.Range("B3").Value = "Bandoleros FCS"
$ws.Range("B4").Value = "cartola scheuer17"
$ws.Range("B5").Value = "seralex"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 71.36
$ws.Range("H2").Value = 34.36
$ws.Range("I2").Value = 37
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 69.56
$ws.Range("H3").Value = 53.06
$ws.Range("I3").Value = 16.5
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 53.06
$ws.Range("H4").Value = 69.56
$ws.Range("I4").Value = -16.5
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 34.36
$ws.Range("H5").Value = 71.36
$ws.Range("I5").Value = -37

# Grupo D
$ws = $wb.Worksheets.Item("Grupo D")
$ws.Range("B2").Value = "Mau Humor F.C."
$ws.Range("B3").Value = "A Lenda Super Vasco F.c"
$ws.Range("B4").Value = "FBC Colorado"
$ws.Range("B5").Value = "Grêmio imortal 36"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 69.76000000000001
$ws.Range("H2").Value = 45.46
$ws.Range("I2").Value = 24.3
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 59.56
$ws.Range("H3").Value = 30.6
$ws.Range("I3").Value = 28.96
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 45.46
$ws.Range("H4").Value = 69.76000000000001
$ws.Range("I4").Value = -24.3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 30.6
$ws.Range("H5").Value = 59.56
$ws.Range("I5").Value = -28.96

# Grupo E
$ws = $wb.Worksheets.Item("Grupo E")
$ws.Range("B2").Value = "Paulo Virgili FC"
$ws.Range("B3").Value = "KillerColorado"
$ws.Range("B4").Value = "Fedato Futebol Clube"
$ws.Range("B5").Value = "FÚRIA LEON"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 73.66
$ws.Range("H2").Value = 50.6
$ws.Range("I2").Value = 23.06
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 57.25
$ws.Range("H3").Value = 49.16
$ws.Range("I3").Value = 8.090000000000003
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 50.6
$ws.Range("H4").Value = 73.66
$ws.Range("I4").Value = -23.06
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 49.16
$ws.Range("H5").Value = 57.25
$ws.Range("I5").Value = -8.090000000000003

# Grupo F
$ws = $wb.Worksheets.Item("Grupo F")
$ws.Range("B2").Value = "lsauer fc"
$ws.Range("B3").Value = "Rolo Compressor ZN"
$ws.Range("B4").Value = "DM Studio"
$ws.Range("B5").Value = "AZURRA82"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 54.36
$ws.Range("H2").Value = 30.06
$ws.Range("I2").Value = 24.3
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 51.05
$ws.Range("H3").Value = 47.2
$ws.Range("I3").Value = 3.849999999999994
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 47.2
$ws.Range("H4").Value = 51.05
$ws.Range("I4").Value = -3.849999999999994
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 30.06
$ws.Range("H5").Value = 54.36
$ws.Range("I5").Value = -24.3

# Grupo G
$ws = $wb.Worksheets.Item("Grupo G")
$ws.Range("B2").Value = "TORRESMO COM PINGA PRO26.1"
$ws.Range("B3").Value = "Tabajara de Inhaua PB1"
$ws.Range("B4").Value = "A Lenda Super Vascão f.c"
$ws.Range("B5").Value = "Grêmio imortal 37"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 69.26000000000001
$ws.Range("H2").Value = 43.1
$ws.Range("I2").Value = 26.16
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 57.76
$ws.Range("H3").Value = 45.3
$ws.Range("I3").Value = 12.46
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 45.3
$ws.Range("H4").Value = 57.76
$ws.Range("I4").Value = -12.46
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 43.1
$ws.Range("H5").Value = 69.26000000000001
$ws.Range("I5").Value = -26.16

# Grupo H
$ws = $wb.Worksheets.Item("Grupo H")
$ws.Range("B2").Value = "Texas Club 2026"
$ws.Range("B3").Value = "Gremiomaniasm"
$ws.Range("B4").Value = "TEAM LOPES 99"
$ws.Range("B5").Value = "Super Vasco f.c"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 59.86
$ws.Range("H2").Value = 50.76
$ws.Range("I2").Value = 9.100000000000001
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 56.65
$ws.Range("H3").Value = 34.76
$ws.Range("I3").Value = 21.89
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 50.76
$ws.Range("H4").Value = 59.86
$ws.Range("I4").Value = -9.100000000000001
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 34.76
$ws.Range("H5").Value = 56.65
$ws.Range("I5").Value = -21.89
